$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3215.7407
$ws.Range("I62").Value = 2940.8823
$ws.Range("J62").Value = 3683
$ws.Range("K62").Value = 2940.8823
$ws.Range("L62").Value = 3683
$ws.Range("M62").Value = -2316.8823
$ws.Range("N62").Value = -4931

$ws.Range("H65").Value = 3215.7407
$ws.Range("I65").Value = 2940.8823
$ws.Range("J65").Value = 3683
$ws.Range("K65").Value = 14704.4115
$ws.Range("L65").Value = 18415
$ws.Range("M65").Value = -11584.4115
$ws.Range("N65").Value = -24655

$ws.Range("H96").Value = 31252344
$ws.Range("I96").Value = 50002164
$ws.Range("J96").Value = 2643
$ws.Range("K96").Value = 150006492
$ws.Range("L96").Value = 7929
$ws.Range("M96").Value = -150005119
$ws.Range("N96").Value = -10675

$ws.Range("H100").Value = 11943019
$ws.Range("I100").Value = 16668231
$ws.Range("K100").Value = 16668231
$ws.Range("M100").Value = -16667690

$ws.Range("H107").Value = 439968.97
$ws.Range("I107").Value = 532346.4399999999
$ws.Range("J107").Value = 1176
$ws.Range("K107").Value = 532346.4399999999
$ws.Range("L107").Value = 1176
$ws.Range("M107").Value = -530426.4399999999
$ws.Range("N107").Value = -5016

$ws.Range("H132").Value = 2438.182
$ws.Range("I132").Value = 2633.7036
$ws.Range("J132").Value = 1558.3334
$ws.Range("K132").Value = 7901.110799999999
$ws.Range("L132").Value = 4675.0002
$ws.Range("M132").Value = -5371.110799999999
$ws.Range("N132").Value = -9735.0002

$ws.Range("H135").Value = 556.6799999999999
$ws.Range("I135").Value = 540.875
$ws.Range("J135").Value = 936
$ws.Range("K135").Value = 4867.875
$ws.Range("L135").Value = 8424
$ws.Range("M135").Value = -2332.875
$ws.Range("N135").Value = -13494

$ws.Range("H137").Value = 971.17145
$ws.Range("I137").Value = 818.44446
$ws.Range("J137").Value = 1486.625
$ws.Range("K137").Value = 2455.33338
$ws.Range("L137").Value = 4459.875
$ws.Range("M137").Value = 94.66661999999997
$ws.Range("N137").Value = -9559.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2337.01
$ws.Range("I32").Value = 2333.989
$ws.Range("K32").Value = 2333.989
$ws.Range("M32").Value = -2046.989

$ws.Range("H74").Value = 3080.0566
$ws.Range("I74").Value = 3223.2292
$ws.Range("J74").Value = 1705.6
$ws.Range("K74").Value = 3223.2292
$ws.Range("L74").Value = 1705.6
$ws.Range("M74").Value = -2349.2292
$ws.Range("N74").Value = -3453.6

$ws.Range("H77").Value = 3080.0566
$ws.Range("I77").Value = 3223.2292
$ws.Range("J77").Value = 1705.6
$ws.Range("K77").Value = 16116.146
$ws.Range("L77").Value = 8528
$ws.Range("M77").Value = -11748.146
$ws.Range("N77").Value = -17264

$ws.Range("H102").Value = 1415.7894
$ws.Range("I102").Value = 1229.4117
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1229.4117
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 392.5882999999999
$ws.Range("N102").Value = -6244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 17533
$ws.Range("I97").Value = 628
$ws.Range("J97").Value = 25985.5
$ws.Range("K97").Value = 628
$ws.Range("L97").Value = 25985.5
$ws.Range("M97").Value = 363
$ws.Range("N97").Value = -27967.5

$ws.Range("H134").Value = 1172.3286
$ws.Range("I134").Value = 826.01697
$ws.Range("J134").Value = 3029.818
$ws.Range("K134").Value = 2478.05091
$ws.Range("L134").Value = 9089.454000000002
$ws.Range("M134").Value = 56.94909000000007
$ws.Range("N134").Value = -14159.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3236.7944
$ws.Range("I31").Value = 3016.5193
$ws.Range("J31").Value = 3782.238
$ws.Range("K31").Value = 3016.5193
$ws.Range("L31").Value = 3782.238
$ws.Range("M31").Value = -2721.5193
$ws.Range("N31").Value = -4372.237999999999

$ws.Range("H34").Value = 3236.7944
$ws.Range("I34").Value = 3016.5193
$ws.Range("J34").Value = 3782.238
$ws.Range("K34").Value = 3016.5193
$ws.Range("L34").Value = 3782.238
$ws.Range("M34").Value = -2814.5193
$ws.Range("N34").Value = -4186.237999999999

$ws.Range("H58").Value = 1051.338
$ws.Range("I58").Value = 580.56366
$ws.Range("J58").Value = 2669.625
$ws.Range("K58").Value = 580.56366
$ws.Range("L58").Value = 2669.625
$ws.Range("M58").Value = -377.56366
$ws.Range("N58").Value = -3075.625

$ws.Range("H132").Value = 2062.45
$ws.Range("I132").Value = 1175
$ws.Range("J132").Value = 5612.25
$ws.Range("K132").Value = 3525
$ws.Range("L132").Value = 16836.75
$ws.Range("M132").Value = -995
$ws.Range("N132").Value = -21896.75

$ws.Range("H136").Value = 1051.338
$ws.Range("I136").Value = 580.56366
$ws.Range("J136").Value = 2669.625
$ws.Range("K136").Value = 1741.69098
$ws.Range("L136").Value = 8008.875
$ws.Range("M136").Value = 808.3090199999999
$ws.Range("N136").Value = -13108.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 67
$ws.Range("I11").Value = 67
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 201
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -61
$ws.Range("N11").ClearContents()

$ws.Range("H80").Value = 4986.6665
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4986.6665
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 14959.9995
$ws.Range("N80").Value = -16831.9995
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 4986.6665
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4986.6665
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 44879.9985
$ws.Range("N83").Value = -54239.9985
$ws.Range("M83").ClearContents()

$ws.Range("H92").Value = 1009.2
$ws.Range("J92").Value = 962
$ws.Range("L92").Value = 2886
$ws.Range("N92").Value = -5382

$ws.Range("H96").Value = 5214
$ws.Range("J96").Value = 5214
$ws.Range("L96").Value = 15642
$ws.Range("N96").Value = -19760

$ws.Range("H118").Value = 1136.5555
$ws.Range("I118").Value = 564.5
$ws.Range("J118").Value = 1300
$ws.Range("K118").Value = 1693.5
$ws.Range("L118").Value = 3900
$ws.Range("M118").Value = -450.5
$ws.Range("N118").Value = -6386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5186.5815
$ws.Range("I70").Value = 4571.8
$ws.Range("J70").Value = 6040.4443
$ws.Range("K70").Value = 4571.8
$ws.Range("L70").Value = 6040.4443
$ws.Range("M70").Value = -4301.8
$ws.Range("N70").Value = -6580.4443

$ws.Range("H73").Value = 5186.5815
$ws.Range("I73").Value = 4571.8
$ws.Range("J73").Value = 6040.4443
$ws.Range("K73").Value = 4571.8
$ws.Range("L73").Value = 6040.4443
$ws.Range("M73").Value = -3635.8
$ws.Range("N73").Value = -7912.4443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5930.5
$ws.Range("I93").Value = 7793.857
$ws.Range("J93").Value = 1582.6666
$ws.Range("K93").Value = 7793.857
$ws.Range("L93").Value = 1582.6666
$ws.Range("M93").Value = -6545.857
$ws.Range("N93").Value = -4078.6666

$ws.Range("H132").Value = 4753.7075
$ws.Range("I132").Value = 4419.353
$ws.Range("J132").Value = 6377.7144
$ws.Range("K132").Value = 13258.059
$ws.Range("L132").Value = 19133.1432
$ws.Range("M132").Value = -10728.059
$ws.Range("N132").Value = -24193.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1803.8125
$ws.Range("I81").Value = 1686.9
$ws.Range("J81").Value = 1998.6666
$ws.Range("K81").Value = 3373.8
$ws.Range("L81").Value = 3997.3332
$ws.Range("M81").Value = -2312.8
$ws.Range("N81").Value = -6119.3332

$ws.Range("H84").Value = 1803.8125
$ws.Range("I84").Value = 1686.9
$ws.Range("J84").Value = 1998.6666
$ws.Range("K84").Value = 16869
$ws.Range("L84").Value = 19986.666
$ws.Range("M84").Value = -11565
$ws.Range("N84").Value = -30594.666

$ws.Range("H136").Value = 5850959
$ws.Range("I136").Value = 7576308
$ws.Range("J136").Value = 11316.154
$ws.Range("K136").Value = 22728924
$ws.Range("L136").Value = 33948.462
$ws.Range("M136").Value = -22726374
$ws.Range("N136").Value = -39048.462
